# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff/Handback DateTime"
# timestamps in the handback status workbook to reflect the new report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for row 2
$overview.Range("G2").Value = "2016-08-19 21:08:18"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$zhcn.Range("H2").Value = "2016-08-19 21:08:14"
$zhcn.Range("K2").Value = "2016-08-19 21:08:30"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
# (Handoff Datetime shares the same underlying timestamp text as Overview's
# Generate Date, so it also becomes 2016-08-19 21:08:18.)
$dede.Range("H2").Value = "2016-08-19 21:08:18"
$dede.Range("K2").Value = "2016-08-19 21:08:36"
